# Auto-generated Excel COM edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 2110.1667
$ws.Cells.Item(9, 9).Value = 1754.7142
$ws.Cells.Item(9, 10).Value = 2607.8
$ws.Cells.Item(9, 11).Value = 1754.7142
$ws.Cells.Item(9, 12).Value = 2607.8
$ws.Cells.Item(9, 13).Value = -1585.7142
$ws.Cells.Item(9, 14).Value = -2945.8
$ws.Cells.Item(15, 8).Value = 1604.2394
$ws.Cells.Item(15, 9).Value = 1604.2394
$ws.Cells.Item(15, 11).Value = 4812.718199999999
$ws.Cells.Item(15, 13).Value = -4643.718199999999
$ws.Cells.Item(33, 8).Value = 404.31818
$ws.Cells.Item(33, 9).Value = 404.8889
$ws.Cells.Item(33, 10).Value = 401.75
$ws.Cells.Item(33, 11).Value = 404.8889
$ws.Cells.Item(33, 12).Value = 401.75
$ws.Cells.Item(33, 13).Value = -175.8889
$ws.Cells.Item(33, 14).Value = -859.75
$ws.Cells.Item(40, 8).Value = 4409.553
$ws.Cells.Item(40, 9).Value = 4080.8572
$ws.Cells.Item(40, 10).Value = 4675.0386
$ws.Cells.Item(40, 11).Value = 4080.8572
$ws.Cells.Item(40, 12).Value = 4675.0386
$ws.Cells.Item(40, 13).Value = -3905.8572
$ws.Cells.Item(40, 14).Value = -5025.0386
$ws.Cells.Item(106, 8).Value = 12148.479
$ws.Cells.Item(106, 9).Value = 7764.5293
$ws.Cells.Item(106, 10).Value = 24569.666
$ws.Cells.Item(106, 11).Value = 7764.5293
$ws.Cells.Item(106, 12).Value = 24569.666
$ws.Cells.Item(106, 13).Value = -7133.5293
$ws.Cells.Item(106, 14).Value = -25831.666
$ws.Cells.Item(132, 8).Value = 1527.1428
$ws.Cells.Item(132, 9).Value = 1579.1621
$ws.Cells.Item(132, 10).Value = 1142.2
$ws.Cells.Item(132, 11).Value = 4737.4863
$ws.Cells.Item(132, 12).Value = 3426.6
$ws.Cells.Item(132, 13).Value = -2207.4863
$ws.Cells.Item(132, 14).Value = -8486.6
$ws.Cells.Item(135, 8).Value = 1282.3793
$ws.Cells.Item(135, 9).Value = 1047.1177
$ws.Cells.Item(135, 10).Value = 2996.4285
$ws.Cells.Item(135, 11).Value = 9424.059300000001
$ws.Cells.Item(135, 12).Value = 26967.8565
$ws.Cells.Item(135, 13).Value = -6889.059300000001
$ws.Cells.Item(135, 14).Value = -32037.8565
$ws.Cells.Item(137, 8).Value = 3074.4363
$ws.Cells.Item(137, 9).Value = 2470.5833
$ws.Cells.Item(137, 11).Value = 7411.749899999999
$ws.Cells.Item(137, 13).Value = -4861.749899999999
$ws.Cells.Item(138, 8).Value = 2518.2654
$ws.Cells.Item(138, 9).Value = 1200.6342
$ws.Cells.Item(138, 10).Value = 3466.0352
$ws.Cells.Item(138, 11).Value = 3601.9026
$ws.Cells.Item(138, 12).Value = 10398.1056
$ws.Cells.Item(138, 13).Value = 1538.0974
$ws.Cells.Item(138, 14).Value = -20678.1056
$ws.Cells.Item(141, 8).Value = 1975.7894
$ws.Cells.Item(141, 9).Value = 1565
$ws.Cells.Item(141, 11).Value = 4695
$ws.Cells.Item(141, 13).Value = 485

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 8841.571
$ws.Cells.Item(2, 9).Value = 1003.6667
$ws.Cells.Item(2, 10).Value = 22949.8
$ws.Cells.Item(2, 11).Value = 1003.6667
$ws.Cells.Item(2, 12).Value = 22949.8
$ws.Cells.Item(2, 13).Value = -890.6667
$ws.Cells.Item(2, 14).Value = -23175.8
$ws.Cells.Item(32, 8).Value = 2375.488
$ws.Cells.Item(32, 9).Value = 2272.111
$ws.Cells.Item(32, 10).Value = 5166.6665
$ws.Cells.Item(32, 11).Value = 2272.111
$ws.Cells.Item(32, 12).Value = 5166.6665
$ws.Cells.Item(32, 13).Value = -1985.111
$ws.Cells.Item(32, 14).Value = -5740.6665
$ws.Cells.Item(45, 8).Value = 100001140
$ws.Cells.Item(45, 9).Value = 111112150
$ws.Cells.Item(45, 11).Value = 111112150
$ws.Cells.Item(45, 13).Value = -111111773
$ws.Cells.Item(88, 8).Value = 4113.9443
$ws.Cells.Item(88, 9).Value = 4207.143
$ws.Cells.Item(88, 10).Value = 4054.6365
$ws.Cells.Item(88, 11).Value = 4207.143
$ws.Cells.Item(88, 12).Value = 4054.6365
$ws.Cells.Item(88, 13).Value = -3801.143
$ws.Cells.Item(88, 14).Value = -4866.636500000001
$ws.Cells.Item(91, 8).Value = 4113.9443
$ws.Cells.Item(91, 9).Value = 4207.143
$ws.Cells.Item(91, 10).Value = 4054.6365
$ws.Cells.Item(91, 11).Value = 4207.143
$ws.Cells.Item(91, 12).Value = 4054.6365
$ws.Cells.Item(91, 13).Value = -2803.143
$ws.Cells.Item(91, 14).Value = -6862.636500000001
$ws.Cells.Item(116, 8).Value = 8841.571
$ws.Cells.Item(116, 9).Value = 1003.6667
$ws.Cells.Item(116, 10).Value = 22949.8
$ws.Cells.Item(116, 11).Value = 1003.6667
$ws.Cells.Item(116, 12).Value = 22949.8
$ws.Cells.Item(116, 13).Value = 1290.3333
$ws.Cells.Item(116, 14).Value = -27537.8
$ws.Cells.Item(132, 8).Value = 4133.8975
$ws.Cells.Item(132, 9).Value = 3438.627
$ws.Cells.Item(132, 10).Value = 8368.727999999999
$ws.Cells.Item(132, 11).Value = 10315.881
$ws.Cells.Item(132, 12).Value = 25106.184
$ws.Cells.Item(132, 13).Value = -7785.880999999999
$ws.Cells.Item(132, 14).Value = -30166.184

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 8841.571
$ws.Cells.Item(3, 9).Value = 1003.6667
$ws.Cells.Item(3, 10).Value = 22949.8
$ws.Cells.Item(3, 11).Value = 1003.6667
$ws.Cells.Item(3, 12).Value = 22949.8
$ws.Cells.Item(3, 13).Value = -889.6667
$ws.Cells.Item(3, 14).Value = -23177.8
$ws.Cells.Item(22, 8).Value = 935
$ws.Cells.Item(22, 9).Value = 608.8
$ws.Cells.Item(22, 10).Value = 1750.5
$ws.Cells.Item(22, 11).Value = 608.8
$ws.Cells.Item(22, 12).Value = 1750.5
$ws.Cells.Item(22, 13).Value = -435.8
$ws.Cells.Item(22, 14).Value = -2096.5
$ws.Cells.Item(132, 8).Value = 59816
$ws.Cells.Item(132, 10).Value = 59816
$ws.Cells.Item(132, 12).Value = 59816
$ws.Cells.Item(132, 14).Value = -69936

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 29610.738
$ws.Cells.Item(31, 9).Value = 3342.6296
$ws.Cells.Item(31, 10).Value = 76893.336
$ws.Cells.Item(31, 11).Value = 3342.6296
$ws.Cells.Item(31, 12).Value = 76893.336
$ws.Cells.Item(31, 13).Value = -3047.6296
$ws.Cells.Item(31, 14).Value = -77483.336
$ws.Cells.Item(34, 8).Value = 29610.738
$ws.Cells.Item(34, 9).Value = 3342.6296
$ws.Cells.Item(34, 10).Value = 76893.336
$ws.Cells.Item(34, 11).Value = 3342.6296
$ws.Cells.Item(34, 12).Value = 76893.336
$ws.Cells.Item(34, 13).Value = -3140.6296
$ws.Cells.Item(34, 14).Value = -77297.336
$ws.Cells.Item(132, 8).Value = 2300.1458
$ws.Cells.Item(132, 9).Value = 2009.4762
$ws.Cells.Item(132, 10).Value = 4334.8335
$ws.Cells.Item(132, 11).Value = 6028.4286
$ws.Cells.Item(132, 12).Value = 13004.5005
$ws.Cells.Item(132, 13).Value = -3498.4286
$ws.Cells.Item(132, 14).Value = -18064.5005
$ws.Cells.Item(134, 8).Value = 1854.9265
$ws.Cells.Item(134, 9).Value = 1161.3617
$ws.Cells.Item(134, 10).Value = 3407.1904
$ws.Cells.Item(134, 11).Value = 3484.0851
$ws.Cells.Item(134, 12).Value = 10221.5712
$ws.Cells.Item(134, 13).Value = -949.0850999999998
$ws.Cells.Item(134, 14).Value = -15291.5712

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 1417.8667
$ws.Cells.Item(92, 9).Value = 457.94116
$ws.Cells.Item(92, 10).Value = 2673.1538
$ws.Cells.Item(92, 11).Value = 1373.82348
$ws.Cells.Item(92, 12).Value = 8019.4614
$ws.Cells.Item(92, 13).Value = -125.82348
$ws.Cells.Item(92, 14).Value = -10515.4614
$ws.Cells.Item(97, 8).Value = 1800
$ws.Cells.Item(97, 10).Value = 1800
$ws.Cells.Item(97, 12).Value = 5400
$ws.Cells.Item(97, 14).Value = -6392
$ws.Cells.Item(116, 8).Value = 1780934.5
$ws.Cells.Item(116, 9).Value = 1335880.4
$ws.Cells.Item(116, 10).Value = 2671043
$ws.Cells.Item(116, 11).Value = 4007641.2
$ws.Cells.Item(116, 12).Value = 8013129
$ws.Cells.Item(116, 13).Value = -4004199.2
$ws.Cells.Item(116, 14).Value = -8020013

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7921.8335
$ws.Cells.Item(80, 9).Value = 2173.3333
$ws.Cells.Item(80, 10).Value = 13670.333
$ws.Cells.Item(80, 11).Value = 2173.3333
$ws.Cells.Item(80, 12).Value = 13670.333
$ws.Cells.Item(80, 13).Value = -1175.3333
$ws.Cells.Item(80, 14).Value = -15666.333
$ws.Cells.Item(83, 8).Value = 7921.8335
$ws.Cells.Item(83, 9).Value = 2173.3333
$ws.Cells.Item(83, 10).Value = 13670.333
$ws.Cells.Item(83, 11).Value = 10866.6665
$ws.Cells.Item(83, 12).Value = 68351.66500000001
$ws.Cells.Item(83, 13).Value = -5874.666499999999
$ws.Cells.Item(83, 14).Value = -78335.66500000001
$ws.Cells.Item(97, 8).Value = 1434
$ws.Cells.Item(97, 9).Value = 1472.8636
$ws.Cells.Item(97, 10).Value = 1291.5
$ws.Cells.Item(97, 11).Value = 1472.8636
$ws.Cells.Item(97, 12).Value = 1291.5
$ws.Cells.Item(97, 13).Value = -976.8635999999999
$ws.Cells.Item(97, 14).Value = -2283.5
$ws.Cells.Item(102, 8).Value = 3411.2104
$ws.Cells.Item(102, 9).Value = 2746.077
$ws.Cells.Item(102, 10).Value = 4852.3335
$ws.Cells.Item(102, 11).Value = 2746.077
$ws.Cells.Item(102, 12).Value = 4852.3335
$ws.Cells.Item(102, 13).Value = -1124.077
$ws.Cells.Item(102, 14).Value = -8096.3335
$ws.Cells.Item(122, 8).Value = 8849
$ws.Cells.Item(122, 9).Value = 9700.521000000001
$ws.Cells.Item(122, 10).Value = 6672.8887
$ws.Cells.Item(122, 11).Value = 29101.563
$ws.Cells.Item(122, 12).Value = 20018.6661
$ws.Cells.Item(122, 13).Value = -26651.563
$ws.Cells.Item(122, 14).Value = -24918.6661
$ws.Cells.Item(132, 8).Value = 3599.6924
$ws.Cells.Item(132, 9).Value = 3304.9546
$ws.Cells.Item(132, 11).Value = 9914.863799999999
$ws.Cells.Item(132, 13).Value = -7384.863799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 5053.561
$ws.Cells.Item(132, 9).Value = 4519.7427
$ws.Cells.Item(132, 10).Value = 8167.5
$ws.Cells.Item(132, 11).Value = 13559.2281
$ws.Cells.Item(132, 12).Value = 24502.5
$ws.Cells.Item(132, 13).Value = -11029.2281
$ws.Cells.Item(132, 14).Value = -29562.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 402.92856
$ws.Cells.Item(113, 9).Value = 294.35715
$ws.Cells.Item(113, 10).Value = 511.5
$ws.Cells.Item(113, 11).Value = 883.0714499999999
$ws.Cells.Item(113, 12).Value = 1534.5
$ws.Cells.Item(113, 13).Value = 1286.92855
$ws.Cells.Item(113, 14).Value = -5874.5
$ws.Cells.Item(122, 8).Value = 3664.2693
$ws.Cells.Item(122, 9).Value = 2527.5264
$ws.Cells.Item(122, 10).Value = 6749.7144
$ws.Cells.Item(122, 11).Value = 7582.5792
$ws.Cells.Item(122, 12).Value = 20249.1432
$ws.Cells.Item(122, 13).Value = -5132.5792
$ws.Cells.Item(122, 14).Value = -25149.1432
$ws.Cells.Item(136, 8).Value = 1622.5698
$ws.Cells.Item(136, 9).Value = 1178.2162
$ws.Cells.Item(136, 10).Value = 4362.75
$ws.Cells.Item(136, 11).Value = 3534.6486
$ws.Cells.Item(136, 12).Value = 13088.25
$ws.Cells.Item(136, 13).Value = -984.6486000000004
$ws.Cells.Item(136, 14).Value = -18188.25
